# Insert a new data row at row 169 (pushing existing rows 169-242 down to 170-243)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 169; this shifts row 169..242 down to 170..243
# and the new blank row 169 inherits formatting (incl. date style on column D) from
# the row it was inserted in front of.
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new record values.
$ws.Cells.Item(169, 1).Value = 11
$ws.Cells.Item(169, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(169, 3).Value = "Bíobío"
$ws.Cells.Item(169, 4).Value = 45141
$ws.Cells.Item(169, 5).Value = 8
$ws.Cells.Item(169, 6).Value = "Fruta"
$ws.Cells.Item(169, 7).Value = 100102
$ws.Cells.Item(169, 8).Value = "Cítricos"
$ws.Cells.Item(169, 9).Value = 100102004
$ws.Cells.Item(169, 10).Value = "Mandarina"
$ws.Cells.Item(169, 11).Value = "Clemenuless"
$ws.Cells.Item(169, 12).Value = "Primera"
$ws.Cells.Item(169, 13).Value = 100
$ws.Cells.Item(169, 14).Value = 7500
$ws.Cells.Item(169, 15).Value = 8000
$ws.Cells.Item(169, 16).Value = 7750
$ws.Cells.Item(169, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(169, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(169, 19).Value = 775
$ws.Cells.Item(169, 20).Value = 10
